$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header: exactMatch (I1)
$ws.Range("I1").Value = "exactMatch"

# Row 9 - exactMatch test row
$ws.Range("A9").Value = 8
$ws.Range("B9").Value = "test exactMatch"
$ws.Range("C9").Value = "test exactMatch"
$ws.Range("D9").Value = "test exactMatch"
$ws.Range("I9").Value = "http://registry.it.csiro.au/sandbox/csiro/utils/commondef/1"

# New header: closeMatch (J1)
$ws.Range("J1").Value = "closeMatch"

# Row 10 - closeMatch test row
$ws.Range("A10").Value = 9
$ws.Range("J10").Value = "http://registry.it.csiro.au/def/environment/feature/water-features"
$ws.Range("B10").Value = "test closeMatch"
$ws.Range("C10").Value = "test closeMatch"
$ws.Range("D10").Value = "test closeMatch"

# Row 11 - related test row
$ws.Range("A11").Value = 10
$ws.Range("B11").Value = "test related"

# New header: related (K1)
$ws.Range("K1").Value = "related"

$ws.Range("K11").Value = "http://registry.it.csiro.au/def/environment/feature/GroundWaterBody"

# Bold header formatting to match existing header style
$ws.Range("I1:K1").Font.Bold = $true

$ws.Range("K11").Select()
